$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.716.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.246.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.30%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'322.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'101.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.14%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.579"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.40%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.554"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.87%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.67%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0830"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.42%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.27%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -2.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.588.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.35%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.854"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.35%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.05%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.244.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.25%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'43.610.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.12%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0985"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.12%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'65.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.34%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'236.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.86%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.69%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.11%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.09%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.21%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'36.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.26%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'160.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.59%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'20.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.76%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0852"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.94%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.05%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +7.80%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -4.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.01%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.48%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.46%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'4.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.47%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'15.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +21.50%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0318"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.72%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.22%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.794.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.03%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'75.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.18%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'82.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.66%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'5.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.53%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'58.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.63%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.23%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'103.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.16%  "
$ws.Range("E51").Style = "Normal"
